# Insert a new data row at row 22 (pushing the existing rows 22-68 down to
# 23-69) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = 1
$ws.Range("B22").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C22").Value = "Arica y Parinacota"
$ws.Range("D22").Value = 44721
$ws.Range("E22").Value = 15
$ws.Range("F22").Value = 100112012
$ws.Range("G22").Value = "Espinaca"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 2750
$ws.Range("N22").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 917
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = "Hortaliza"
